# Applies the Xhosa-translation copy edits described by the diff.
# Each Find/Replace targets an exact source string (MatchCase=$true,
# MatchWholeWord=$false, MatchWildcards=$false) so unrelated partial
# matches elsewhere in the document are not touched. Where a source
# string is itself a substring of another (e.g. the three
# "Wamkelekile kwi-ParentText" variants, or "Thatha ikhefu" vs.
# "...ukuThatha ikhefu..."), the most specific / first-in-document
# variant is replaced first with wdReplaceOne(=1) so only the intended
# occurrence is changed before the broader pattern is searched again.

$d = $word.ActiveDocument
$wdReplaceOne = 1
$wdReplaceAll = 2

function Replace-Text($find, $replace, $mode) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, $mode) | Out-Null
}

# 1. "Wamkelekile kwi-ParentText! " (first occurrence, has trailing "! ")
Replace-Text "Wamkelekile kwi-ParentText! " "Wamkelekile kwi ParentText! " $wdReplaceOne

# 2. Robot intro paragraph
Replace-Text "Ndingu ______, umkhokheli wakho. Nangona kubonakala ngathi nddingumntu, ndilirobhothi elenziwe yi-Parenting for Lifelong Health no-UNICEF ukuzokuxhasa kuhambo lwakho lokuba ngumzali. " "Ndingu ______, umkhokheli wakho. Nangona kubonakala ngathi ndingumntu, ndiyirobhothi eyenziwe yi-Parenting for Lifelong Health no-UNICEF uzokuxhasa kuhambo lwakho lokuba ngumzali. " $wdReplaceAll

# 3. "Wamkelekile kwi-ParentText!" (second occurrence, ends with "!" no trailing space)
Replace-Text "Wamkelekile kwi-ParentText!" "Wamkelekile kwi ParentText!" $wdReplaceOne

# 4. "Make sijonge ukuba isebenza njani i-ParentText. "
Replace-Text "Make sijonge ukuba isebenza njani i-ParentText. " "Masikhe sijonge ukuba isebenza kanjani i-ParentText. " $wdReplaceAll

# 5. "Ngoku, make sijonge ukuba sibonakala njani isisfundo se-ParentText. "
Replace-Text "Ngoku, make sijonge ukuba sibonakala njani isisfundo se-ParentText. " "Ngoku, makhe sijonge ukuba sibonakala njani isifundo se-ParentText. " $wdReplaceAll

# 6. Reminder-message paragraph
Replace-Text "Uzakufumana umyalezo okukhumbuzayo ukuba ugqibe isifundo sakho. Ukuba awuwubonanga umyalezo kulungile! Ungabuyela kwi ParentText nanini na ukuzikhumbuza ngesifundo sakho. Kwaye ukuba ikuphosile, kulungile! Usenokubuyela kwi ParentText nangaliphi na ixesha ukuze uqhubele phambili nesifundo sakho." "Uzakufumana umyalezo yonke imihla okukhumbuzayo ukuba ugqibe isifundo sakho. Kwaye ukuba sikuphosile, kulungile! Usenokubuyela kwi ParentText nangaliphi na ixesha ukuze uqhubekele phambili nesifundo sakho." $wdReplaceAll

# 7. "Isifundo ngasinye siquka imibuzo, ..."
Replace-Text "Isifundo ngasinye siquka imibuzo, imifanekiso, amacebiso nemidlalo emnandi ukuze uzame ukuwenza ekhaya nomntwana okanye nosapho lwakho." "Isifundo ngasinye siquka imibuzo, imifanekiso, Iingcebiso kunye nemidlalo emnandi ukuze uzame ukuyenza ekhaya nomntwana okanye nosapho lwakho." $wdReplaceAll

# 8. Replaced entirely with a screenshot placeholder instruction
Replace-Text "Xa uziva ukuba ufuna uncedo. bhala u-Menu okanye Uncedo ekugqibeleni kwesifundo ukuze ufumane inkxaso eyongezelelweyo" "{Show screenshots of the quiz, tips, comics, and home activity}" $wdReplaceAll

# 9. "Xa ufuna uncedo, bhala MENU okanye HELP ..."
Replace-Text "Xa ufuna uncedo, bhala MENU okanye HELP ekupheleni kwesifundo sakho ukufumana inkxaso eyongezelelweyo. " "Ukuba ukhe waxinga okanye ufuna uncedo, bhala MENU okanye NCEDA ekupheleni kwesifundo sakho uzokufumana inkxaso eyongezelelweyo. " $wdReplaceAll

# 10. "Xa ubhala UNCEDO nanini na, ..." (long help/safety paragraph, split over 3 sentences)
Replace-Text "Xa ubhala UNCEDO nanini na, ungafumana ulwazi ngoovimba abakhoyo ekuhlaleni ukumelana nobundlobongela bentsapho, ubundlobongela ngokwesondo, impilo yengqondo, nezinye iimeko zongxamiseko. " "Xa ubhala NCEDA nanini na, ungafumana ulwazi ngezixhobo ezikhoyo ekuhlaleni ukumelana nobundlobongela bosapho, ubundlobongela ngokwesondo, impilo yengqondo, okanye nezinye iimeko zongxamiseko. " $wdReplaceAll
Replace-Text "Ulwazi lwakho lukhuselekile apha: Alukho ulwazi okuzokwabelwana ngalo ngaphandle kwemvume yakho okanye luthengiswe ukwenza inzuzo. Umyalezo owuthumelayo uvaliwe kwaye utshixelwe kwiseva ekhuselekileyo. Le miyalezo oyithumelayo inoguqulelo oluntsokothileyo kwaye itshixelwe kwiseva ekhuselekileyo. " "Iinkcukacha zakho zikhuselekile apha: Akukho nanye ekuzokwabelwana ngayo ngaphandle kwemvume yakho kwaye azizukuthengiswa ukwenza inzuzo. Le miyalezo oyithumelayo inoguqulelo oluntsokothileyo kwaye itshixelwe kwiseva ekhuselekileyo. " $wdReplaceAll
Replace-Text "Khumbula, xa kukho umntu okwazi ukufikelela kwifoni yakho xa ingatshixwanga angakwazi ukubona imiyalezo yakho. Xa ngamanye amaxesha uthumela ulwazi olunobuzaza kwaye loonto ikukhathaza cima yonke imilayezo kwifoni yakho. Ngoko ke, ukuba uthumela ulwazi ulunobuzaza kwaye unexhala, cima imiyalezo kwifowuni yakho. " "Khumbula, nabani na okwaziyo ukufikelela kwifoni yakho xa ingatshixwanga angakwazi ukubona imiyalezo yakho. Ngoko ke, ukuba uthumela ulwazi olunobuzaza kwaye unexhala, cima imiyalezo kwifowuni yakho. " $wdReplaceAll

# 11. "Ukuba khona kwakho apha kubonisa ..."
Replace-Text "Ukuba khona kwakho apha kubonisa ukuba ukukhathalele kangakanani ukubonelela inkxaso yomntwana wakho. " "Ukuba khona kwakho apha kubonisa ukuba ukukhathalele kangakanani umntwana wakho eyona nkxaso. " $wdReplaceAll

# 12. "Yilonto uyenza nomntwana ezakwenza umehluko. "
Replace-Text "Yilonto uyenza nomntwana ezakwenza umehluko. " "Yinto oyenzayo nomntwana ezakwenza umahluko. " $wdReplaceAll

# 13. "I-ParentText iya kubonelela ngamacebiso ..."
Replace-Text "I-ParentText iya kubonelela ngamacebiso ngezifundo ezizakunceda ngobudlelwane bakho nomntwana wakho. Kukuwe ukusebenzisa lamacebiso uzame uwaprakthize. Kuxhomekeke kuwe ukuba uwasebenzise la macebiso!" "I-ParentText iya kubonelela ngeengcebiso ngezifundo ezizakunceda kubudlelwane bakho nomntwana wakho. Kuxhomekeke kuwe ukuba uyazisebenzisa ezingcebiso!" $wdReplaceAll

# 14. "Enkosi kakhulu ukumamela! ..."
Replace-Text "Enkosi kakhulu ukumamela! Ungayifumana le-vidiyo nanini na xa usiya kwi-Menu. Siyathemba ukonwabele ukuba kwi-ParentText nokuthi uzolusebenzisa ulwazi olufumene apha! Ungafikelela kulevidiyo nangaliphi na ixesha ngeMENYU. Siyathemba uya kukonwabela ukusebenzisa i ParentText kwaye wenze lukhulu kuyo! " "Enkosi kakhulu ngokumamela! Ungafikelela kulevidiyo nangaliphi na ixesha ngeMENYU. Siyathemba uya kukonwabela ukusebenzisa i ParentText kwaye wenze lukhulu kuyo! " $wdReplaceAll

# 15. "Wamkelekile kwi-ParentText" (third, bare occurrence - last remaining after edits 1 & 3)
Replace-Text "Wamkelekile kwi-ParentText" "Wamkelekile kwi ParentText" $wdReplaceOne

# 16. "Molo! Uziva njani njengangoku? Unayo imizuzwana eyi 30?"
Replace-Text "Molo! Uziva njani njengangoku? Unayo imizuzwana eyi 30?" "Molo! Uziva njani njengangoku? Unayo imizuzwana engamashumi amathathu?" $wdReplaceAll

# 17. "Ngaphambi kokuba siqale kwi-ParentText, masithathe ikhefu elikhawulezileyo."
Replace-Text "Ngaphambi kokuba siqale kwi-ParentText, masithathe ikhefu elikhawulezileyo." "Ngaphambi kokuba uqale kwi ParentText, masithi nqumama ngokukhawuleza kunye." $wdReplaceAll

# 18. "Thatha ikhefu" (standalone; first occurrence in the doc, before the longer sentences below)
Replace-Text "Thatha ikhefu" "Thatha Ikhefu" $wdReplaceOne

# 19. "Phefumla nzulu."
Replace-Text "Phefumla nzulu." "Phefumla nzuuulu." $wdReplaceAll

# 20. "Wuve umphefumlo wakho ungena, uphuma emzimbeni wakho."
Replace-Text "Wuve umphefumlo wakho ungena, uphuma emzimbeni wakho." "Wuve umoya ungena, uphuma, emzimbeni wakho." $wdReplaceAll

# 21. "Qwalasela ukuba uluva kweyiphi indawo emzimbeni wakho uxinzelelo."
Replace-Text "Qwalasela ukuba uluva kweyiphi indawo emzimbeni wakho uxinzelelo." "Qwalasela ukuba uluva kweyiphi indawo uxinzelelo emzimbeni wakho." $wdReplaceAll

# 22. "Zama ukuyolula okanye ukuyiphumza londayo."
Replace-Text "Zama ukuyolula okanye ukuyiphumza londayo." "Zama ukuyiphumza londawo." $wdReplaceAll

# 23. "Xa sele ukulungele ukuwavula amehlo, ungawavula."
Replace-Text "Xa sele ukulungele ukuwavula amehlo, ungawavula." "Xa sele ulungile, vula amehlo kwakhona." $wdReplaceAll

# 24. "Ngoku qwalasela ukuba ngabe uziva ngokwehlukileyo kunangokuya"
Replace-Text "Ngoku qwalasela ukuba ngabe uziva ngokwehlukileyo kunangokuya" "Ngoku, qwalasela ukuba ingaba uziva ngokwahlukileyo kunangokuya" $wdReplaceAll

# 25. "ubusaqala ukuthatha eli khefu."
Replace-Text "ubusaqala ukuthatha eli khefu." "xa ubusaqala lomsebenzi." $wdReplaceAll

# 26. "Zama ukuThatha ikhefu nanini na xa uziva unomsindo, unoxinzelelo okanye ukhathazekile."
Replace-Text "Zama ukuThatha ikhefu nanini na xa uziva unomsindo, unoxinzelelo okanye ukhathazekile." "Zama ukuThatha Ikhefu nanini na xa uziva unomsindo, unoxinzelelo, okanye ukhathazekile." $wdReplaceAll

# 27. "Nokuphefumla kambalwa okanye ukunxulumana nomhlaba, kungenza umehluko."
Replace-Text "Nokuphefumla kambalwa okanye ukunxulumana nomhlaba, kungenza umehluko." "Nokuphefumla nzulu kambalwa, okanye ukunxulumana nomhlaba, kungenza umehluko." $wdReplaceAll

# 28. "Kwaye ungayiThatha ikhefu nabantwana wakho omncinci okanye omdala!"
Replace-Text "Kwaye ungayiThatha ikhefu nabantwana wakho omncinci okanye omdala!" "UngayiThatha Ikhefu nomntwana wakho omncinci okanye omdala!" $wdReplaceAll

Write-Output "Done applying Xhosa copy edits."
